# Anonymize "fedcore" -> "approach" in the header rows, and give the
# "fedcore"/"approach" and "change" header cells (columns C/D, and F/G on
# the computational_comparison sheet) a top+bottom border (matching the
# style already used for the merged "original" header cell group), plus
# clear the stray empty inline-string cell left in G5.

$xlNone = -4142   # xlLineStyleNone
$xlThin = 1       # xlContinuous / thin weight for LineStyle

function Set-TopBottomBorder($range) {
    # Reset to the plain/default cell format first (these header cells end
    # up on the default, non-bold style with just a border applied), then
    # add a top + bottom thin border, no left/right border (= border id 4).
    $range.ClearFormats()
    $range.Borders.Item(8).LineStyle = $xlThin    # xlEdgeTop
    $range.Borders.Item(9).LineStyle = $xlThin    # xlEdgeBottom
}

function Set-TopBottomRightBorder($range) {
    # top + bottom + right thin border, no left border (= border id 5)
    $range.ClearFormats()
    $range.Borders.Item(8).LineStyle = $xlThin    # xlEdgeTop
    $range.Borders.Item(9).LineStyle = $xlThin    # xlEdgeBottom
    $range.Borders.Item(10).LineStyle = $xlThin   # xlEdgeRight
}

$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-TopBottomBorder $ws1.Range("C1")
Set-TopBottomRightBorder $ws1.Range("D1")

$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-TopBottomBorder $ws2.Range("C1")
Set-TopBottomRightBorder $ws2.Range("D1")
Set-TopBottomBorder $ws2.Range("F1")
Set-TopBottomRightBorder $ws2.Range("G1")

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell at G5
$ws2.Range("G5").ClearContents()
